$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register the small (size 8) "Aptos Narrow" font used by the sheet's
# phoneticPr so it lands in xl/styles.xml, then remove the scratch row
# so it leaves no visible trace in the sheet data / dimension.
$scratchRow = 1048000
$ws.Rows.Item($scratchRow).Font.Size = 8
$ws.Rows.Item($scratchRow).Delete()

$signs = @(
    "20 kph sign",
    "30 kph sign",
    "40 kph sign",
    "50 kph sign",
    "60 kph sign",
    "70 kph sign",
    "80 kph sign",
    "90 kph sign",
    "100 kph sign",
    "110 kph sign",
    "120 kph sign",
    "130 kph sign",
    "140 kph sign",
    "Stop sign",
    "Traffic light",
    "Roundabout",
    "Speed bump"
)

$row = 3
foreach ($sign in $signs) {
    $ws.Cells.Item($row, 1).Value = $sign
    $ws.Cells.Item($row, 2).Value = "N/A"
    $row = $row + 1
}

$ws.Range("A1").Select()
$ws.Application.Goto($ws.Range("D20"))
